# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "Periodo Mora" of the existing worker row (was 2507) to the new period 2508
$ws.Range("E16").Value = "2508"

# The second "Periodo Mora" row (2506) is no longer part of this statement -
# remove that entire worksheet row, shifting everything below it up.
$ws.Rows("17:17").Delete()

# Recalculate the totals now that only one period remains:
# VALOR MORA total (was sum of the two periods, 113880) -> single period value
$ws.Range("E11").Value = 56940

# Cant. Periodos (count of periods) (was 2) -> 1
$ws.Range("F13").Value = 1
